$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B (2) through AD (30) hold the match data; column A is just the
# running index and stays put on each row.
$firstCol = 2
$lastCol = 30

function Get-RowValues($row) {
    $vals = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $vals += , ($ws.Cells.Item($row, $c).Value2)
    }
    return $vals
}

function Set-RowValues($row, $vals) {
    for ($i = 0; $i -lt $vals.Count; $i++) {
        $ws.Cells.Item($row, $firstCol + $i).Value2 = $vals[$i]
    }
}

# --- Rows 36 and 37: the two fixtures were swapped ---
$row36 = Get-RowValues 36
$row37 = Get-RowValues 37
Set-RowValues 36 $row37
Set-RowValues 37 $row36

# --- Rows 189, 190, 191: the three fixtures were cyclically rotated ---
# new189 = old191, new190 = old189, new191 = old190
$row189 = Get-RowValues 189
$row190 = Get-RowValues 190
$row191 = Get-RowValues 191

Set-RowValues 189 $row191
Set-RowValues 190 $row189
Set-RowValues 191 $row190
